# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 9 (pushing the existing rows 9-23
# down to 10-24), matching the updated "Hortaliza, Terminal Hortofrutícola
# Agro Chillán - Berenjena" consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:23 down one row, creating a fresh (empty) row 9.
$ws.Rows(9).Insert()

# Populate the new row 9 with the latest weekly record.
$ws.Cells.Item(9, 1).Value  = 7
$ws.Cells.Item(9, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value  = "Ñuble"
$ws.Cells.Item(9, 4).Value  = 44624
$ws.Cells.Item(9, 5).Value  = 16
$ws.Cells.Item(9, 6).Value  = 100112001
$ws.Cells.Item(9, 7).Value  = "Berenjena"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 10000
$ws.Cells.Item(9, 12).Value = 11000
$ws.Cells.Item(9, 13).Value = 10500
$ws.Cells.Item(9, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(9, 15).Value = "Región Metropolitana"
$ws.Cells.Item(9, 16).Value = 175
$ws.Cells.Item(9, 17).Value = 60
$ws.Cells.Item(9, 18).Value = "Hortaliza"
